# Refresh cryptos list snapshot (prices / 1h volume %) from upstream source
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.462.34'
$ws.Range("E2").Value = '  +8.51%  '

$ws.Range("D3").Value = '2.579.88'
$ws.Range("E3").Value = '  +10.10%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'504.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.22%  '

$ws.Range("D6").Value = "'157.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.27%  '

$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +25.47%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '2.582.49'
$ws.Range("E9").Value = '  +10.03%  '

$ws.Range("D10").Value = "'6.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +15.40%  '

$ws.Range("E11").Value = '  +6.97%  '

$ws.Range("E12").Value = '  +7.29%  '

$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").Value = '3.029.87'
$ws.Range("E14").Value = '  +9.98%  '

$ws.Range("D15").Value = '59.374.22'
$ws.Range("E15").Value = '  +7.53%  '

$ws.Range("D16").Value = "'21.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.22%  '

$ws.Range("E17").Value = '  +5.86%  '

$ws.Range("D18").Value = '2.593.97'
$ws.Range("E18").Value = '  +10.49%  '

$ws.Range("D19").Value = "'4.74"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'334.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.77%  '

$ws.Range("E21").Value = '  +7.89%  '

$ws.Range("E22").Value = '  +8.59%  '

$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").Value = "'60.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.89%  '

$ws.Range("E25").Value = '  +5.93%  '

$ws.Range("E26").Value = '  +9.87%  '

$ws.Range("D27").Value = '2.701.26'
$ws.Range("E27").Value = '  +10.02%  '

$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("D29").Value = '0.0₃0823'
$ws.Range("E29").Value = '  +9.98%  '

$ws.Range("E30").Value = '  +3.37%  '

$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("D32").Value = "'157.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.03%  '

$ws.Range("D33").Value = "'19.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.46%  '

$ws.Range("E34").Value = '  +7.41%  '

$ws.Range("E35").Value = '  +10.00%  '

$ws.Range("D36").Value = "'3.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.51%  '

$ws.Range("E37").Value = '  +9.69%  '

$ws.Range("D38").Value = "'0.849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.99%  '

$ws.Range("E39").Value = '  +12.51%  '

$ws.Range("E40").Value = '  +9.04%  '

$ws.Range("D41").Value = "'35.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.98%  '

$ws.Range("D42").Value = "'289.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.54%  '

$ws.Range("D43").Value = "'0.102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.04%  '

$ws.Range("E44").Value = '  +8.92%  '

$ws.Range("D45").Value = "'0.0564"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.19%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").Value = "'19.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +15.54%  '

$ws.Range("E48").Value = '  +6.34%  '

$ws.Range("E49").Value = '  +16.26%  '

$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = "'10.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = "'4.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.25%  '
